$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "71.143.30", "  +2.97%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.837.73", "  +1.72%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.04%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "706.71", "  +12.30%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "172.84", "  +4.13%  "),
    @(7, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.837.45", "  +1.76%  "),
    @(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.999", "  -0.02%  "),
    @(9, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.528", "  +1.28%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.163", "  +2.46%  "),
    @(11, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "7.42", "  +9.24%  "),
    @(12, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.465", "  +1.05%  "),
    @(13, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000248", "  +1.50%  "),
    @(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "36.27", "  +3.53%  "),
    @(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.477.06", "  +1.68%  "),
    @(16, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.832.43", "  +1.15%  "),
    @(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "71.105.75", "  +2.93%  "),
    @(18, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "17.91", "  +1.84%  "),
    @(19, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "7.24", "  +3.28%  "),
    @(20, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.115", "  +0.73%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "11.49", "  +19.96%  "),
    @(22, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "481.85", "  +3.66%  "),
    @(23, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.719", "  +1.61%  "),
    @(24, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "84.04", "  +1.26%  "),
    @(25, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0000144", "  -0.90%  "),
    @(26, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "12.38", "  +3.08%  "),
    @(27, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "10.41", "  +3.94%  "),
    @(28, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "2.15", "  -0.54%  "),
    @(29, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.982.97", "  +1.61%  "),
    @(30, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.21", "  +20.26%  "),
    @(31, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  -0.05%  "),
    @(32, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.34", "  +2.81%  "),
    @(33, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "7.53", "  +5.92%  "),
    @(34, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "29.86", "  +4.61%  "),
    @(35, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.180", "  +2.67%  "),
    @(36, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "9.26", "  +3.01%  "),
    @(37, "RenzoRestakedETH", "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth", "3.782.68", "  +1.64%  "),
    @(38, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  +0.00%  "),
    @(39, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.103", "  +1.39%  "),
    @(40, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "3.45", "  +4.31%  "),
    @(41, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "6.02", "  +3.55%  "),
    @(42, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "2.25", "  +16.46%  "),
    @(43, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.973", "  +0.89%  "),
    @(44, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  -0.11%  "),
    @(45, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  +0.02%  "),
    @(46, "FLOKI", "https://coinranking.com/coin/fmHk13Rqw+floki-floki", "0.000308", "  +13.29%  "),
    @(47, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "49.81", "  +6.46%  "),
    @(48, "ONDO", "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo", "1.45", "  +1.29%  "),
    @(49, "Arweave", "https://coinranking.com/coin/7XWg41D1+arweave-ar", "45.51", "  +4.83%  "),
    @(50, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "160.29", "  +1.60%  "),
    @(51, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.302", "  +2.03%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]
    $bCell.Style = "Normal"

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $row[2]
    $cCell.Style = "Normal"

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]
    $eCell.Style = "Normal"
}
